$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.213.95"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "1.855.43"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("D4").Value = "'0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("B5").Value = "XRP"
$ws.Range("C5").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.7009"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'241.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("D8").Value = "'0.3091"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.58%  "

$ws.Range("D9").Value = "'0.07726"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").Value = "'23.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "

$ws.Range("E11").Value = "  -2.40%  "

$ws.Range("D12").Value = "1.864.22"
$ws.Range("E12").Value = "  -0.67%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.098"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.55%  "

$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'92.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").Value = "'0.6876"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "'6.499"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.33%  "

$ws.Range("D17").Value = "'0.000008414"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "29.221.89"
$ws.Range("E18").Value = "  -0.58%  "

$ws.Range("D19").Value = "'249.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.07%  "

$ws.Range("D20").Value = "2.115.29"
$ws.Range("E20").Value = "  -1.68%  "

$ws.Range("D21").Value = "'12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").Value = "'7.526"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "'0.9999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "

$ws.Range("D25").Value = "'0.1519"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "

$ws.Range("D26").Value = "'160.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").Value = "'8.854"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.58%  "

$ws.Range("E28").Value = "  -1.09%  "

$ws.Range("D29").Value = "'1.559"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.02%  "

$ws.Range("D30").Value = "'4.232"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").Value = "'4.199"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("D32").Value = "'1.191"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("D33").Value = "'0.05184"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.43%  "

$ws.Range("D34").Value = "'0.7634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.38%  "

$ws.Range("D35").Value = "'1.844"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.57%  "

$ws.Range("D36").Value = "'1.163"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.58%  "

$ws.Range("D37").Value = "'2.711"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").Value = "'0.01861"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.12%  "

$ws.Range("D39").Value = "1.221.45"
$ws.Range("E39").Value = "  -1.96%  "

$ws.Range("D40").Value = "'2.726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.68%  "

$ws.Range("D41").Value = "'0.8967"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.48%  "

$ws.Range("D42").Value = "'109.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "

$ws.Range("D43").Value = "'0.9994"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").Value = "'5.531"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.50%  "

$ws.Range("D45").Value = "2.012.41"
$ws.Range("E45").Value = "  -2.67%  "

$ws.Range("E46").Value = "  -4.26%  "

$ws.Range("D47").Value = "'65.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.42%  "

$ws.Range("D48").Value = "'0.5178"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").Value = "'9.536"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.59%  "

$ws.Range("D50").Value = "'1.750"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").Value = "'7.015"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
